$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('I8').Value = 'sd'
$ws.Range('J8').Value = 'Statement-non-opinion'
$ws.Range('I24').Value = 'ba'
$ws.Range('J24').Value = 'Appreciation'
$ws.Range('I63').Value = 'sd'
$ws.Range('J63').Value = 'Statement-non-opinion'
$ws.Range('I66').Value = 'sv'
$ws.Range('J66').Value = 'Statement-opinion'
$ws.Range('I67').Value = 'aa'
$ws.Range('J67').Value = 'Agree/Accept'
$ws.Range('I72').Value = 'b'
$ws.Range('J72').Value = 'Acknowledge (Backchannel)'
$ws.Range('I92').Value = 'ba'
$ws.Range('J92').Value = 'Appreciation'
$ws.Range('I101').Value = 'b'
$ws.Range('J101').Value = 'Acknowledge (Backchannel)'
$ws.Range('I104').Value = 'ba'
$ws.Range('J104').Value = 'Appreciation'
$ws.Range('I106').Value = 'sv'
$ws.Range('J106').Value = 'Statement-opinion'
$ws.Range('I120').Value = 'qy'
$ws.Range('J120').Value = 'Yes-No-Question'
$ws.Range('I125').Value = 'sv'
$ws.Range('J125').Value = 'Statement-opinion'
$ws.Range('I127').Value = 'sd'
$ws.Range('J127').Value = 'Statement-non-opinion'
$ws.Range('I133').Value = 'sd'
$ws.Range('J133').Value = 'Statement-non-opinion'
$ws.Range('I134').Value = 'sv'
$ws.Range('J134').Value = 'Statement-opinion'
$ws.Range('I147').Value = 'sd'
$ws.Range('J147').Value = 'Statement-non-opinion'
$ws.Range('I149').Value = 'sd'
$ws.Range('J149').Value = 'Statement-non-opinion'
$ws.Range('I183').Value = 'ba'
$ws.Range('J183').Value = 'Appreciation'
$ws.Range('I187').Value = 'sd'
$ws.Range('J187').Value = 'Statement-non-opinion'
$ws.Range('I191').Value = 'sd'
$ws.Range('J191').Value = 'Statement-non-opinion'
$ws.Range('I193').Value = 'sd'
$ws.Range('J193').Value = 'Statement-non-opinion'
$ws.Range('I204').Value = 'ba'
$ws.Range('J204').Value = 'Appreciation'
$ws.Range('I210').Value = 'sd'
$ws.Range('J210').Value = 'Statement-non-opinion'
$ws.Range('I215').Value = 'sv'
$ws.Range('J215').Value = 'Statement-opinion'
$ws.Range('I216').Value = '%'
$ws.Range('J216').Value = 'Uninterpretable'
$ws.Range('I232').Value = 'sd'
$ws.Range('J232').Value = 'Statement-non-opinion'
$ws.Range('I241').Value = 'sd'
$ws.Range('J241').Value = 'Statement-non-opinion'
$ws.Range('I248').Value = 'sv'
$ws.Range('J248').Value = 'Statement-opinion'
$ws.Range('I251').Value = 'ba'
$ws.Range('J251').Value = 'Appreciation'
$ws.Range('I254').Value = 'b'
$ws.Range('J254').Value = 'Acknowledge (Backchannel)'
$ws.Range('I255').Value = 'aa'
$ws.Range('J255').Value = 'Agree/Accept'
$ws.Range('I256').Value = 'aa'
$ws.Range('J256').Value = 'Agree/Accept'
$ws.Range('I257').Value = 'sd'
$ws.Range('J257').Value = 'Statement-non-opinion'
$ws.Range('I261').Value = '%'
$ws.Range('J261').Value = 'Uninterpretable'
$ws.Range('I285').Value = 'sv'
$ws.Range('J285').Value = 'Statement-opinion'
$ws.Range('I292').Value = 'aa'
$ws.Range('J292').Value = 'Agree/Accept'
$ws.Range('I295').Value = 'sd'
$ws.Range('J295').Value = 'Statement-non-opinion'
$ws.Range('I296').Value = 'sd'
$ws.Range('J296').Value = 'Statement-non-opinion'
$ws.Range('I323').Value = 'sd'
$ws.Range('J323').Value = 'Statement-non-opinion'
$ws.Range('I326').Value = 'ba'
$ws.Range('J326').Value = 'Appreciation'
$ws.Range('I329').Value = 'sd'
$ws.Range('J329').Value = 'Statement-non-opinion'
$ws.Range('I339').Value = 'ba'
$ws.Range('J339').Value = 'Appreciation'
$ws.Range('I350').Value = 'aa'
$ws.Range('J350').Value = 'Agree/Accept'
$ws.Range('I351').Value = 'sd'
$ws.Range('J351').Value = 'Statement-non-opinion'
$ws.Range('I368').Value = 'b'
$ws.Range('J368').Value = 'Acknowledge (Backchannel)'
$ws.Range('I370').Value = 'b'
$ws.Range('J370').Value = 'Acknowledge (Backchannel)'
$ws.Range('I371').Value = 'ba'
$ws.Range('J371').Value = 'Appreciation'
$ws.Range('I397').Value = 'sd'
$ws.Range('J397').Value = 'Statement-non-opinion'
$ws.Range('I404').Value = 'sv'
$ws.Range('J404').Value = 'Statement-opinion'
$ws.Range('I405').Value = 'sd'
$ws.Range('J405').Value = 'Statement-non-opinion'
$ws.Range('I409').Value = 'sd'
$ws.Range('J409').Value = 'Statement-non-opinion'
$ws.Range('I425').Value = '%'
$ws.Range('J425').Value = 'Uninterpretable'
$ws.Range('I430').Value = 'sv'
$ws.Range('J430').Value = 'Statement-opinion'
$ws.Range('I433').Value = 'sd'
$ws.Range('J433').Value = 'Statement-non-opinion'
$ws.Range('I435').Value = 'sv'
$ws.Range('J435').Value = 'Statement-opinion'
$ws.Range('I443').Value = 'sv'
$ws.Range('J443').Value = 'Statement-opinion'
$ws.Range('I464').Value = 'sd'
$ws.Range('J464').Value = 'Statement-non-opinion'
$ws.Range('I476').Value = 'sd'
$ws.Range('J476').Value = 'Statement-non-opinion'
$ws.Range('I478').Value = 'sv'
$ws.Range('J478').Value = 'Statement-opinion'
$ws.Range('I509').Value = 'ba'
$ws.Range('J509').Value = 'Appreciation'
$ws.Range('I519').Value = 'sd'
$ws.Range('J519').Value = 'Statement-non-opinion'
$ws.Range('I535').Value = 'ba'
$ws.Range('J535').Value = 'Appreciation'
$ws.Range('I544').Value = 'ba'
$ws.Range('J544').Value = 'Appreciation'
$ws.Range('I554').Value = 'sd'
$ws.Range('J554').Value = 'Statement-non-opinion'
$ws.Range('I559').Value = 'sv'
$ws.Range('J559').Value = 'Statement-opinion'
$ws.Range('I561').Value = 'ba'
$ws.Range('J561').Value = 'Appreciation'
$ws.Range('I565').Value = 'ba'
$ws.Range('J565').Value = 'Appreciation'
$ws.Range('I574').Value = 'sd'
$ws.Range('J574').Value = 'Statement-non-opinion'
$ws.Range('I580').Value = 'ba'
$ws.Range('J580').Value = 'Appreciation'
